$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4:H4").NumberFormat = "@"
$ws.Range("A4").Value = "12"
$ws.Range("B4").Value = "34"
$ws.Range("C4").Value = "12"
$ws.Range("D4").Value = "6"
$ws.Range("E4").Value = "6"
$ws.Range("F4").Value = "7"
$ws.Range("G4").Value = "8"
$ws.Range("H4").Value = "9"
